$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Filter roads by type in assessment:
# - Remove the stray space in the "3m/2m/1m" road-type labels
# - Mark the "2 m Weg" / "1 m Weg" rows as "to keep" = yes
$ws.Range("B12").Value = "3m Strasse"
$ws.Range("B16").Value = "2m Weg"
$ws.Range("B17").Value = "1m Weg"
$ws.Range("C16").Value = "yes"
$ws.Range("C17").Value = "yes"

# Row heights settle back to the sheet's normal (auto-fit) values once the
# text is tidied up.
for ($r = 8; $r -le 27; $r++) {
    $ws.Rows($r).RowHeight = 18.75
}
$ws.Rows(10).RowHeight = 44.25
$ws.Rows(11).RowHeight = 44.25
$ws.Rows(21).RowHeight = 44.25
$ws.Rows(12).RowHeight = 69.75
